# Adding a few lines of config for notifications
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 19:20 (pushes the existing "gap" row + the two
# trailing contracts_finances_<region>/contracts_<region> rows down by two)
$ws.Rows("19:20").Insert()

# Row 19: notifications / oos:id
$ws.Range("A19").Value = "notifications"
$ws.Range("B19").Value = "notifications_<region>"
$ws.Range("C19").Value = "/*/d1:notificationEF/oos:id"
$ws.Range("D19").Value = "oos:id"
$ws.Range("G19").Value = "NotificationID"
$ws.Range("E19").Value = '"13768"'
$ws.Range("F19").Value = "integer"

# Row 20: notifications / oos:notificationNumber
$ws.Range("A20").Value = "notifications"
$ws.Range("B20").Value = "notifications_<region>"
$ws.Range("C20").Value = "/*/d1:notificationEF/oos:notificationNumber"
$ws.Range("D20").Value = "oos:notificationNumber"
$ws.Range("G20").Value = "NotificationNumber"
$ws.Range("E20").Value = '"0176100001811000002"'
$ws.Range("F20").Value = "character"

# Match the author's final selection
$ws.Range("A20").Select()
